# Weekly data refresh: insert a new observation as the first data row
# (row 7), pushing the existing rows down by one. The new row carries the
# same market/product metadata as the row it displaces, but with updated
# date, volume, price and origin figures for the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7; Excel shifts rows 7:63 down to 8:64 and the
# used range / dimension grows to A1:T64 automatically.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new week's record.
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 44490
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100108
$ws.Cells.Item(7, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(7, 9).Value = 100108002
$ws.Cells.Item(7, 10).Value = "Mango"
$ws.Cells.Item(7, 11).Value = "Sin especificar"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 30
$ws.Cells.Item(7, 14).Value = 7500
$ws.Cells.Item(7, 15).Value = 8000
$ws.Cells.Item(7, 16).Value = 7750
$ws.Cells.Item(7, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(7, 18).Value = "Perú"
$ws.Cells.Item(7, 19).Value = 1938
$ws.Cells.Item(7, 20).Value = 4
